$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new rows of time-tracking data
$ws.Range("A2").Copy($ws.Range("A6"))
$ws.Range("A6").Value = 44315
$ws.Range("B6").Value = 5.25
$ws.Range("C6").Value = "Сверстал страницы, добавил валидацию полей для создания заявки."

$ws.Range("A2").Copy($ws.Range("A7"))
$ws.Range("A7").Value = 44316
$ws.Range("B7").Value = 1.5
$ws.Range("C7").Value = "Закончил валидацию полей на странице создания заявки."

# Recalculate the total formula in F2 (already =SUM(B:B), just force recalc)
$ws.Range("F2").Formula = "=SUM(B:B)"
$excel.Calculate()

# Update the selected cell to match the new active selection
$ws.Range("B9").Select()
